$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B, C, D, E, G across rows 2-6 (F column unchanged)
$data = @{
    2 = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 3.223369029078222;   E = 0.5333859586016987; G = 8.656069925401464 }
    3 = @{ B = 1.445647641019636;  C = 1.626987699542094;  D = 0.1496068669990043;  E = 0.5333859586016987; G = 3.755628166162433 }
    4 = @{ B = 0.04172184405617529; C = 0.04103571897497393; D = 0.1496068669990043; E = 0.5333859586016987; G = 0.7657503886318522 }
    5 = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 0.1496068669990043;  E = 0.5333859586016987; G = 5.582307763322248 }
    6 = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 0.7210945179870265;  E = 13.86384647080068;  G = 19.48425592650926 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
